# Apply crypto price/volume updates and re-ranking swaps
# Generated from the OOXML diff describing the commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "80.568.63"
$ws.Range("E2").Value = "  +4.81%  "

# Row 3
$ws.Range("D3").Value = "3.193.12"
$ws.Range("E3").Value = "  +2.32%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "639.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.17%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.290"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +28.42%  "

# Row 8
$ws.Range("E8").Value = "  -0.11%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.97%  "

# Row 10
$ws.Range("D10").Value = "3.191.42"
$ws.Range("E10").Value = "  +2.32%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.592"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +16.20%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000266"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +26.63%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.165"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.94%  "

# Row 14
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.64%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.779.75"
$ws.Range("E15").Value = "  +2.28%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.34%  "

# Row 17
$ws.Range("D17").Value = "80.489.70"
$ws.Range("E17").Value = "  +4.77%  "

# Row 18
$ws.Range("D18").Value = "3.190.74"
$ws.Range("E18").Value = "  +2.22%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.71%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +14.20%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "446.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +12.31%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.51%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.44%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.63%  "

# Row 25
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.13%  "

# Row 26
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "3.359.90"
$ws.Range("E26").Value = "  +2.32%  "

# Row 27
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "77.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.74%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.03%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.06%  "

# Row 30
$ws.Range("E30").Value = "  +10.88%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.82%  "

# Row 32
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.996"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "570.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +12.05%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.35%  "

# Row 35
$ws.Range("E35").Value = "  +13.17%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.22%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "22.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.80%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.122"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +19.83%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.413"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.12%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.70%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.68%  "

# Row 43
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "159.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.25%  "

# Row 44
$ws.Range("E44").Value = "  +0.03%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "189.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.70%  "

# Row 46
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.11%  "

# Row 47
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.44%  "

# Row 48
$ws.Range("B48").Value = "ImmutableX"
$ws.Range("C48").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.18%  "

# Row 49
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.775"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.81%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.25%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.70%  "

